$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 122 (pushes existing rows 122-133 down to 123-134)
$ws.Rows.Item(122).Insert()

# Populate the newly inserted row 122 with a new weekly data point
# (same market/variety as the last existing row, one week later than the last date)
$ws.Cells.Item(122, 1).Value = 5
$ws.Cells.Item(122, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(122, 3).Value = "Maule"
$ws.Cells.Item(122, 4).Value = 44918
$ws.Cells.Item(122, 5).Value = 7
$ws.Cells.Item(122, 6).Value = 100112022
$ws.Cells.Item(122, 7).Value = "Arveja Verde"
$ws.Cells.Item(122, 8).Value = "Sin especificar"
$ws.Cells.Item(122, 9).Value = "Primera"
$ws.Cells.Item(122, 10).Value = 500
$ws.Cells.Item(122, 11).Value = 20000
$ws.Cells.Item(122, 12).Value = 20000
$ws.Cells.Item(122, 13).Value = 20000
$ws.Cells.Item(122, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(122, 15).Value = "Carahue"
$ws.Cells.Item(122, 16).Value = 800
$ws.Cells.Item(122, 17).Value = 25
$ws.Cells.Item(122, 18).Value = "Hortaliza"
